$d = $word.ActiveDocument

# Locate the paragraph that contains the M2Doc field text "{m:''.availableTableStyles()}"
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*availableTableStyles()}*") {
        $target = $p.Range
        break
    }
}

$start = $target.Start
$text = $target.Text

# Offsets (relative to the paragraph start) of the two split points:
#  - between "{" and "m"            -> splits run "{m" into "{" + "m"
#  - between "()" and "}"           -> splits run "()}" into "()" + "}"
$splitBraceM = $start + $text.IndexOf("{m") + 1
$closeParenPos = $start + $text.IndexOf("()}") + 2

# --- Split "()}" into "()" and "}" (the "}" run ends up with no rPr at all) ---
$d.Range($closeParenPos, $closeParenPos).InsertParagraphAfter()
$d.Range($closeParenPos, $closeParenPos + 1).Delete()
$d.Range($closeParenPos, $closeParenPos + 1).Delete()
$sel = $word.Selection
$sel.SetRange($closeParenPos, $closeParenPos)
$sel.TypeText("}")

# --- Split "{m" into "{" and "m" (both keep the same <w:lang> rPr) ---
$d.Range($splitBraceM, $splitBraceM).InsertParagraphAfter()
$d.Range($splitBraceM, $splitBraceM + 1).Delete()

Write-Output $d.Content.Text
